$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values: force text storage (NumberFormat "@") so the
# numeric-looking strings (e.g. "42.609.93", "0.999", "1.00") are kept
# verbatim as text instead of being auto-coerced to numbers, then reset
# the cell style back to Normal so no stray numeric format sticks around.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '42.609.93'
$ws.Range("E2").Value = '  -0.76%  '
Set-TextValue "D3" '2.540.54'
$ws.Range("E3").Value = '  -0.37%  '
Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue "D5" '313.56'
$ws.Range("E5").Value = '  +2.89%  '
Set-TextValue "D6" '95.17'
Set-TextValue "D7" '0.577'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +0.03%  '
Set-TextValue "D9" '0.537'
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("E10").Value = '  -1.83%  '
Set-TextValue "D11" '0.0816'
$ws.Range("E11").Value = '  -1.99%  '
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("E13").Value = '  -1.07%  '
Set-TextValue "D14" '2.929.77'
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("E15").Value = '  +3.79%  '
Set-TextValue "D16" '2.533.74'
$ws.Range("E16").Value = '  -0.77%  '
Set-TextValue "D17" '0.861'
$ws.Range("E17").Value = '  -1.58%  '
Set-TextValue "D18" '42.647.01'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("E19").Value = '  -3.07%  '
Set-TextValue "D20" '6.67'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("E21").Value = '  -2.52%  '
Set-TextValue "D22" '71.12'
$ws.Range("E22").Value = '  -1.34%  '
Set-TextValue "D23" '254.75'
$ws.Range("E23").Value = '  -1.03%  '
Set-TextValue "D24" '2.95'
$ws.Range("E24").Value = '  -0.26%  '
Set-TextValue "D26" '27.46'
$ws.Range("E26").Value = '  -2.70%  '
Set-TextValue "D27" '1.00'
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +5.22%  '
Set-TextValue "D29" '39.92'
$ws.Range("E29").Value = '  +4.79%  '
Set-TextValue "D30" '10.07'
$ws.Range("E30").Value = '  -1.57%  '
Set-TextValue "D31" '5.97'
$ws.Range("E31").Value = '  -2.97%  '
Set-TextValue "D32" '155.61'
$ws.Range("E32").Value = '  -1.72%  '
Set-TextValue "D33" '19.68'
$ws.Range("E33").Value = '  +0.52%  '
Set-TextValue "D34" '3.40'
$ws.Range("E34").Value = '  +2.19%  '
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("E36").Value = '  -0.86%  '
Set-TextValue "D37" '2.62'
$ws.Range("E37").Value = '  -0.19%  '
Set-TextValue "D38" '0.112'
$ws.Range("E38").Value = '  -4.12%  '
Set-TextValue "D39" '24.80'
$ws.Range("E39").Value = '  -4.14%  '
$ws.Range("E40").Value = '  -0.49%  '
Set-TextValue "D41" '2.19'
$ws.Range("E41").Value = '  +5.25%  '
Set-TextValue "D42" '3.39'
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("E44").Value = '  -1.00%  '
Set-TextValue "D45" '0.999'
$ws.Range("E45").Value = '  +0.06%  '
Set-TextValue "D46" '2.045.46'
$ws.Range("E46").Value = '  -2.64%  '
$ws.Range("E47").Value = '  -3.89%  '
$ws.Range("E48").Value = '  -0.34%  '
Set-TextValue "D49" '75.66'
$ws.Range("E49").Value = '  +0.84%  '
Set-TextValue "D50" '2.784.50'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("E51").Value = '  -0.24%  '
